$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Name" column header, reusing the same header style (yellow fill)
# already applied to the "Code Group" / "Country" headers in A1/B1.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "Name"

# Fill in the new "Name" column values for the existing rows.
$ws.Range("C2").Value = "Lebanon Gov"
$ws.Range("C3").Value = "Syria 1"

# Size the new column to fit its (longer) content, similar to Excel's
# column auto-fit behavior when a new column is added.
$ws.Columns("C").ColumnWidth = 12.86

# Leave the selection on the cell right below the newly entered data,
# matching where the cursor would land after typing the last value.
[void]$ws.Range("C4").Select()
